$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.620.96'
$ws.Range('E2').Value = '  +1.29%  '

$ws.Range('D3').Value = '1.888.81'
$ws.Range('E3').Value = '  +1.60%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.00'
$ws.Range('E5').Value = '  +0.75%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4827'
$ws.Range('E7').Value = '  +0.75%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2871'
$ws.Range('E8').Value = '  +2.46%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06552'
$ws.Range('E9').Value = '  +1.54%  '

$ws.Range('D10').Value = '1.807.89'
$ws.Range('E10').Value = '  -2.77%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07448'
$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.74'
$ws.Range('E12').Value = '  +2.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.094'
$ws.Range('E13').Value = '  +0.20%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.74'
$ws.Range('E14').Value = '  +0.72%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6667'
$ws.Range('E15').Value = '  +3.29%  '

$ws.Range('D16').Value = '30.594.24'
$ws.Range('E16').Value = '  +1.41%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.21'
$ws.Range('E17').Value = '  +0.36%  '

$ws.Range('E18').Value = '  +0.16%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007578'
$ws.Range('E19').Value = '  +0.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '231.10'
$ws.Range('E20').Value = '  +3.18%  '

$ws.Range('D21').Value = '2.140.36'
$ws.Range('E21').Value = '  +1.93%  '

$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.266'
$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.190'
$ws.Range('E24').Value = '  +1.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.383'
$ws.Range('E25').Value = '  +1.88%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.14'
$ws.Range('E26').Value = '  +2.79%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.68'
$ws.Range('E27').Value = '  +1.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.954'
$ws.Range('E28').Value = '  +1.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1022'
$ws.Range('E29').Value = '  +10.87%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.400'
$ws.Range('E30').Value = '  -2.31%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.327'
$ws.Range('E31').Value = '  +2.18%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.020'
$ws.Range('E32').Value = '  +1.66%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05048'
$ws.Range('E33').Value = '  +1.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.205'
$ws.Range('E34').Value = '  +5.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7497'
$ws.Range('E35').Value = '  +3.49%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +0.50%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.713'
$ws.Range('E37').Value = '  +0.95%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01880'
$ws.Range('E38').Value = '  +2.55%  '

$ws.Range('E39').Value = '  +1.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9204'
$ws.Range('E40').Value = '  +2.28%  '

$ws.Range('E41').Value = '  +1.00%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.05'
$ws.Range('E42').Value = '  +0.80%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4285'
$ws.Range('E43').Value = '  +0.81%  '

$ws.Range('E45').Value = '  -3.95%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.424'
$ws.Range('E46').Value = '  +2.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '64.10'
$ws.Range('E47').Value = '  +0.50%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1277'
$ws.Range('E48').Value = '  -2.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.484'
$ws.Range('E49').Value = '  -0.90%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.995'
$ws.Range('E50').Value = '  +3.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.02'
$ws.Range('E51').Value = '  +0.78%  '
